$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; H=258.90410958904107; I="28/11/2021"; J="Sunday"; K="12" }
    @{ Row=3; H=21.575342465753423; I="28/11/2021"; J="Sunday"; K="1" }
    @{ Row=4; H=73.28767123287672; I="28/11/2021"; J="Sunday"; K="6" }
    @{ Row=5; H=43.15068493150685; I="28/11/2021"; J="Sunday"; K="6" }
    @{ Row=6; H=50.68493150684932; I="28/11/2021"; J="Sunday"; K="12" }
    @{ Row=7; H=0; I="28/11/2021"; J="Sunday"; K="" }
    @{ Row=8; H=9.589041095890412; I="28/11/2021"; J="Sunday"; K="12" }
    @{ Row=9; H=4.10958904109589; I="28/11/2021"; J="Sunday"; K="12" }
    @{ Row=10; H=129.45205479452054; I="29/11/2021"; J="Monday"; K="6" }
    @{ Row=11; H=21.575342465753423; I="29/11/2021"; J="Monday"; K="1" }
    @{ Row=12; H=73.28767123287672; I="29/11/2021"; J="Monday"; K="6" }
    @{ Row=13; H=43.15068493150685; I="29/11/2021"; J="Monday"; K="6" }
    @{ Row=14; H=12.67123287671233; I="29/11/2021"; J="Monday"; K="3" }
    @{ Row=15; H=2.28310502283105; I="29/11/2021"; J="Monday"; K="1" }
    @{ Row=16; H=0; I="29/11/2021"; J="Monday"; K="" }
    @{ Row=17; H=4.10958904109589; I="29/11/2021"; J="Monday"; K="12" }
    @{ Row=18; H=129.45205479452054; I="1/12/2021"; J="Tuesday"; K="6" }
    @{ Row=19; H=21.575342465753423; I="1/12/2021"; J="Tuesday"; K="1" }
    @{ Row=20; H=85.50228310502284; I="1/12/2021"; J="Tuesday"; K="7" }
    @{ Row=21; H=50.34246575342466; I="1/12/2021"; J="Tuesday"; K="7" }
    @{ Row=22; H=25.34246575342466; I="1/12/2021"; J="Tuesday"; K="6" }
    @{ Row=23; H=4.5662100456621; I="1/12/2021"; J="Tuesday"; K="2" }
    @{ Row=24; H=9.589041095890412; I="1/12/2021"; J="Tuesday"; K="12" }
    @{ Row=25; H=8.21917808219178; I="1/12/2021"; J="Tuesday"; K="24" }
    @{ Row=26; H=151.02739726027397; I="2/12/2021"; J="Wednesday"; K="7" }
    @{ Row=27; H=21.575342465753423; I="2/12/2021"; J="Wednesday"; K="1" }
    @{ Row=28; H=146.57534246575344; I="2/12/2021"; J="Wednesday"; K="12" }
    @{ Row=29; H=50.34246575342466; I="2/12/2021"; J="Wednesday"; K="7" }
    @{ Row=30; H=25.34246575342466; I="2/12/2021"; J="Wednesday"; K="6" }
    @{ Row=31; H=2.28310502283105; I="2/12/2021"; J="Wednesday"; K="1" }
    @{ Row=32; H=9.589041095890412; I="2/12/2021"; J="Wednesday"; K="12" }
    @{ Row=33; H=8.21917808219178; I="2/12/2021"; J="Wednesday"; K="24" }
    @{ Row=34; H=129.45205479452054; I="3/12/2021"; J="Thursday"; K="6" }
    @{ Row=35; H=21.575342465753423; I="3/12/2021"; J="Thursday"; K="1" }
    @{ Row=36; H=146.57534246575344; I="3/12/2021"; J="Thursday"; K="12" }
    @{ Row=37; H=50.34246575342466; I="3/12/2021"; J="Thursday"; K="7" }
    @{ Row=38; H=25.34246575342466; I="3/12/2021"; J="Thursday"; K="6" }
    @{ Row=39; H=2.28310502283105; I="3/12/2021"; J="Thursday"; K="1" }
    @{ Row=40; H=9.589041095890412; I="3/12/2021"; J="Thursday"; K="12" }
    @{ Row=41; H=8.21917808219178; I="3/12/2021"; J="Thursday"; K="24" }
    @{ Row=42; H=151.02739726027397; I="4/12/2021"; J="Friday"; K="7" }
    @{ Row=43; H=21.575342465753423; I="4/12/2021"; J="Friday"; K="1" }
    @{ Row=44; H=73.28767123287672; I="4/12/2021"; J="Friday"; K="6" }
    @{ Row=45; H=50.34246575342466; I="4/12/2021"; J="Friday"; K="7" }
    @{ Row=46; H=25.34246575342466; I="4/12/2021"; J="Friday"; K="6" }
    @{ Row=47; H=2.28310502283105; I="4/12/2021"; J="Friday"; K="1" }
    @{ Row=48; H=9.589041095890412; I="4/12/2021"; J="Friday"; K="12" }
    @{ Row=49; H=8.21917808219178; I="4/12/2021"; J="Friday"; K="24" }
    @{ Row=50; H=151.02739726027397; I="5/12/2021"; J="Saturday"; K="7" }
    @{ Row=51; H=21.575342465753423; I="5/12/2021"; J="Saturday"; K="1" }
    @{ Row=52; H=146.57534246575344; I="5/12/2021"; J="Saturday"; K="12" }
    @{ Row=53; H=50.34246575342466; I="5/12/2021"; J="Saturday"; K="7" }
    @{ Row=54; H=25.34246575342466; I="5/12/2021"; J="Saturday"; K="6" }
    @{ Row=55; H=2.28310502283105; I="5/12/2021"; J="Saturday"; K="1" }
    @{ Row=56; H=9.589041095890412; I="5/12/2021"; J="Saturday"; K="12" }
    @{ Row=57; H=8.21917808219178; I="5/12/2021"; J="Saturday"; K="24" }
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 8).Value = $item.H

    $ws.Cells.Item($item.Row, 9).NumberFormat = "@"
    $ws.Cells.Item($item.Row, 9).Value = $item.I

    $ws.Cells.Item($item.Row, 10).NumberFormat = "@"
    $ws.Cells.Item($item.Row, 10).Value = $item.J

    $ws.Cells.Item($item.Row, 11).NumberFormat = "@"
    $ws.Cells.Item($item.Row, 11).Value = $item.K
}

Write-Output "done"